$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct passives form factor (0805 -> 0402) ---
# R2, R3, R4, R5, R6, RPROG: RESISTOR device/package (rows 6-11)
foreach ($r in 6..11) {
    $ws.Range("C$r").Value = "R-EU_R0402"
    $ws.Range("D$r").Value = "R0402"
}

# C2: CAPACITOR device/package
$ws.Range("C2").Value = "C-EUC0402"
$ws.Range("D2").Value = "C0402"

# --- Apply left/center alignment to the updated Device/Package cells ---
# Build the combined alignment format on an out-of-the-way scratch cell first
# (setting HorizontalAlignment then VerticalAlignment directly on the target
# range would otherwise leave a stray "horizontal-only" intermediate style
# behind in the workbook's style table), then stamp that finished style onto
# the real target range in a single assignment, and finally wipe the scratch
# cell so it doesn't affect the sheet's used range.
$scratch = $ws.Range("Z100")
$scratch.HorizontalAlignment = -4131
$scratch.VerticalAlignment = -4108
$fmtRange = $ws.Range("C2:D2,C6:D11")
$fmtRange.Style = $scratch.Style
$scratch.Clear()

# --- Update selection to reflect the edited cells ---
$ws.Range("C2:D2").Select()
